# Applies cryptos list update (prices / volume%) per commit
# "Updated cryptos list on Sat Jul 13 22:36:16 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.485.99"
$ws.Range("E2").Value = "  +1.39%  "

$ws.Range("D3").Value = "'3.152.36"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'528.60"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("D6").Value = "'139.47"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  +10.43%  "

$ws.Range("D10").Value = "'0.439"
$ws.Range("E10").Value = "  +6.35%  "

$ws.Range("E11").Value = "  +3.91%  "

$ws.Range("E12").Value = "  +2.34%  "

$ws.Range("D13").Value = "'3.691.60"
$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").Value = "'25.61"
$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("D16").Value = "'58.485.08"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.21"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.143.64"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").Value = "'12.91"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").Value = "'8.10"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "'371.97"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").Value = "'0.528"
$ws.Range("E23").Value = "  +4.48%  "

$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").Value = "'8.35"
$ws.Range("E27").Value = "  +14.51%  "

$ws.Range("D28").Value = "'0.0₃0852"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("D29").Value = "'22.30"
$ws.Range("E29").Value = "  +4.43%  "

$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").Value = "'5.12"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'6.30"
$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("D35").Value = "'156.63"
$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("E36").Value = "  +2.63%  "

$ws.Range("D37").Value = "'2.698.48"
$ws.Range("E37").Value = "  +6.39%  "

$ws.Range("D38").Value = "'24.90"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0690"
$ws.Range("E39").Value = "  +2.88%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.68"
$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("D41").Value = "'4.26"
$ws.Range("E41").Value = "  +5.43%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0292"
$ws.Range("E42").Value = "  +8.75%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.720"
$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("D44").Value = "'39.01"
$ws.Range("E44").Value = "  +3.22%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").Value = "'3.191.37"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("E47").Value = "  +11.47%  "

$ws.Range("D48").Value = "'6.19"
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").Value = "'20.00"
$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("D51").Value = "'0.745"
$ws.Range("E51").Value = "  +0.96%  "
